$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "key"
$ws.Range("B1").Value = "this.firstName"

# Update row 2: keep A2 as RegisterWithoutLastName, clear B2
$ws.Range("A2").Value = "RegisterWithoutLastName"
$ws.Range("B2").ClearContents()

# Update selection to C3
$ws.Range("C3").Select()
